# Auto-generated edit script: applies the numeric corrections from the
# 'chore: update Sheets via scheduled runner' commit to each sheet's
# profit-calculation columns (H:N).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 3187.45
$ws.Range("I98").Value = 2014.5
$ws.Range("J98").Value = 4946.875
$ws.Range("K98").Value = 2014.5
$ws.Range("L98").Value = 4946.875
$ws.Range("M98").Value = -516.5
$ws.Range("N98").Value = -7942.875
$ws.Range("H113").Value = 5383
$ws.Range("I113").Value = 4645
$ws.Range("J113").Value = 6490
$ws.Range("K113").Value = 4645
$ws.Range("L113").Value = 6490
$ws.Range("M113").Value = -1391
$ws.Range("H122").Value = 3187.45
$ws.Range("I122").Value = 2014.5
$ws.Range("J122").Value = 4946.875
$ws.Range("K122").Value = 6043.5
$ws.Range("L122").Value = 14840.625
$ws.Range("M122").Value = -3593.5
$ws.Range("N122").Value = -19740.625
$ws.Range("H138").Value = 4509.951
$ws.Range("I138").Value = 2580.682
$ws.Range("J138").Value = 6743.8423
$ws.Range("K138").Value = 7742.045999999999
$ws.Range("L138").Value = 20231.5269
$ws.Range("M138").Value = -2602.045999999999
$ws.Range("N138").Value = -30511.5269

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 6221
$ws.Range("I31").Value = 6221
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 6221
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -5927
$ws.Range("H32").Value = 3865.0537
$ws.Range("I32").Value = 2632.7563
$ws.Range("J32").Value = 10273
$ws.Range("K32").Value = 2632.7563
$ws.Range("L32").Value = 10273
$ws.Range("M32").Value = -2345.7563
$ws.Range("H63").Value = 3784.95
$ws.Range("I63").Value = 2058.25
$ws.Range("J63").Value = 6375
$ws.Range("K63").Value = 2058.25
$ws.Range("L63").Value = 6375
$ws.Range("M63").Value = -1372.25
$ws.Range("H66").Value = 3784.95
$ws.Range("I66").Value = 2058.25
$ws.Range("J66").Value = 6375
$ws.Range("K66").Value = 10291.25
$ws.Range("L66").Value = 31875
$ws.Range("M66").Value = -6859.25
$ws.Range("H103").Value = 25092.834
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 25092.834
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 25092.834
$ws.Range("N103").Value = -27436.834

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 633
$ws.Range("I94").Value = 633
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 633
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -182
$ws.Range("N94").ClearContents()
$ws.Range("H105").Value = 1656.381
$ws.Range("I105").Value = 1443.3334
$ws.Range("J105").Value = 2189
$ws.Range("K105").Value = 1443.3334
$ws.Range("L105").Value = 2189
$ws.Range("M105").Value = 303.6666
$ws.Range("N105").Value = -5683
$ws.Range("H132").Value = 45290
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 45290
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 45290
$ws.Range("N132").Value = -55410
$ws.Range("H134").Value = 4636.107
$ws.Range("I134").Value = 5367.4
$ws.Range("J134").Value = 3792.3076
$ws.Range("K134").Value = 16102.2
$ws.Range("L134").Value = 11376.9228
$ws.Range("M134").Value = -13567.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1472700.6
$ws.Range("I31").Value = 2501334.5
$ws.Range("J31").Value = 3223.6072
$ws.Range("K31").Value = 2501334.5
$ws.Range("L31").Value = 3223.6072
$ws.Range("M31").Value = -2501039.5
$ws.Range("N31").Value = -3813.6072
$ws.Range("H34").Value = 1472700.6
$ws.Range("I34").Value = 2501334.5
$ws.Range("J34").Value = 3223.6072
$ws.Range("K34").Value = 2501334.5
$ws.Range("L34").Value = 3223.6072
$ws.Range("M34").Value = -2501132.5
$ws.Range("N34").Value = -3627.6072
$ws.Range("H132").Value = 3920.2222
$ws.Range("I132").Value = 2434.25
$ws.Range("J132").Value = 6081.636
$ws.Range("K132").Value = 7302.75
$ws.Range("L132").Value = 18244.908
$ws.Range("M132").Value = -4772.75
$ws.Range("N132").Value = -23304.908

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 827.2787
$ws.Range("I107").Value = 524.8387
$ws.Range("J107").Value = 1139.8
$ws.Range("K107").Value = 1574.5161
$ws.Range("L107").Value = 3419.4
$ws.Range("M107").Value = 345.4838999999999
$ws.Range("H121").Value = 48947.715
$ws.Range("I121").Value = 566.6667
$ws.Range("J121").Value = 85233.5
$ws.Range("K121").Value = 1700.0001
$ws.Range("L121").Value = 255700.5
$ws.Range("M121").Value = -390.0001
$ws.Range("N121").Value = -258320.5
$ws.Range("H131").Value = 1646.75
$ws.Range("I131").Value = 2869.0908
$ws.Range("J131").Value = 1108.92
$ws.Range("K131").Value = 8607.2724
$ws.Range("L131").Value = 3326.76
$ws.Range("M131").Value = -3567.2724
$ws.Range("N131").Value = -13406.76
$ws.Range("H132").Value = 4139
$ws.Range("I132").Value = 2673.75
$ws.Range("J132").Value = 10000
$ws.Range("K132").Value = 24063.75
$ws.Range("L132").Value = 90000
$ws.Range("M132").Value = -21533.75
$ws.Range("N132").Value = -95060
$ws.Range("H133").Value = 6343.3335
$ws.Range("I133").Value = 6343.3335
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 19030.0005
$ws.Range("L133").Value = 0
$ws.Range("M133").Value = -13970.0005
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1172.0358
$ws.Range("I97").Value = 898.3333
$ws.Range("J97").Value = 1993.1428
$ws.Range("K97").Value = 898.3333
$ws.Range("L97").Value = 1993.1428
$ws.Range("M97").Value = -402.3333
$ws.Range("H122").Value = 9191.6
$ws.Range("I122").Value = 11633.333
$ws.Range("J122").Value = 8145.143
$ws.Range("K122").Value = 34899.999
$ws.Range("L122").Value = 24435.429
$ws.Range("M122").Value = -32449.999
$ws.Range("N122").Value = -29335.429

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 55558450
$ws.Range("I16").Value = 83335090
$ws.Range("J16").Value = 5165.6665
$ws.Range("K16").Value = 83335090
$ws.Range("L16").Value = 5165.6665
$ws.Range("M16").Value = -83334920
$ws.Range("N16").Value = -5505.6665
$ws.Range("H22").Value = 250003000
$ws.Range("I22").Value = 333334000
$ws.Range("J22").Value = 10000
$ws.Range("K22").Value = 333334000
$ws.Range("L22").Value = 10000
$ws.Range("M22").Value = -333333705
$ws.Range("H27").Value = 250003000
$ws.Range("I27").Value = 333334000
$ws.Range("J27").Value = 10000
$ws.Range("K27").Value = 333334000
$ws.Range("L27").Value = 10000
$ws.Range("M27").Value = -333333893
$ws.Range("H40").Value = 2708.0908
$ws.Range("I40").Value = 1958.6
$ws.Range("J40").Value = 3332.6667
$ws.Range("K40").Value = 1958.6
$ws.Range("L40").Value = 3332.6667
$ws.Range("M40").Value = -1822.6
$ws.Range("H82").Value = 2840.9565
$ws.Range("I82").Value = 2065
$ws.Range("J82").Value = 3687.4546
$ws.Range("K82").Value = 2065
$ws.Range("L82").Value = 3687.4546
$ws.Range("M82").Value = -1704
$ws.Range("N82").Value = -4409.4546
$ws.Range("H85").Value = 2840.9565
$ws.Range("I85").Value = 2065
$ws.Range("J85").Value = 3687.4546
$ws.Range("K85").Value = 2065
$ws.Range("L85").Value = 3687.4546
$ws.Range("M85").Value = -817
$ws.Range("N85").Value = -6183.4546
$ws.Range("H122").Value = 2930.7437
$ws.Range("I122").Value = 2526.6333
$ws.Range("J122").Value = 4277.778
$ws.Range("K122").Value = 7579.8999
$ws.Range("L122").Value = 12833.334
$ws.Range("M122").Value = -5129.8999
$ws.Range("N122").Value = -17733.334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 4168852.5
$ws.Range("I126").Value = 1996.5714
$ws.Range("J126").Value = 10002451
$ws.Range("K126").Value = 5989.7142
$ws.Range("L126").Value = 30007353
$ws.Range("M126").Value = -3519.7142
